$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update raw data values in row 5 (G5, H5, I5)
$ws.Range("G5").Value = 21
$ws.Range("H5").Value = 22
$ws.Range("I5").Value = 4

# Fix total row (row 15) summation formulas to include row 2 through 14
# instead of row 5 through 14
$ws.Range("B15").Formula = "=SUM(B2:B14)"
$ws.Range("C15").Formula = "=SUM(C2:C14)"
$ws.Range("E15").Formula = "=SUM(E2:E14)"
$ws.Range("G15").Formula = "=SUM(G2:G14)"
$ws.Range("H15").Formula = "=SUM(H2:H14)"
$ws.Range("I15").Formula = "=SUM(I2:I14)"

# Update the active selection to reflect the last-edited cell
$ws.Range("J5").Select()
